# Apply the "scrapper_Service error handling" edit:
#  - A1 text changes from "3273112" to "3273114" (kept as text, bold + wrap-text style)
#  - a new row 2 is added with A2 = "6SL32105BE211UV0" using the same bold + wrap-text style
#  - B1's old numeric value (5) is removed, but its original (pre-edit) style/look is retained
#  - row heights / selection / default column width are adjusted to match the new layout

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give A1 a bold, wrapped-text look (this introduces the new font/style used by A1 & A2) ---
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").WrapText = $true

# Replace A1's content with the new value "3273114".
# (Using a formula + paste-as-values round trip keeps the result a genuine text/shared-string
#  cell instead of Excel auto-converting the numeric-looking text into a number.)
$ws.Range("A1").Formula = '="3273114"'
$ws.Range("A1").Copy()
$ws.Range("A1").PasteSpecial(-4163)  # xlPasteValues

# --- Add the new second row with the new part number, reusing A1's formatting ---
$ws.Range("A2").Value = "6SL32105BE211UV0"
$ws.Range("A1").Copy()
$ws.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- B1 loses its value but keeps its former appearance ---
$ws.Range("B1").ClearContents()

# --- Row heights to accommodate the wrapped / bold text ---
$ws.Rows.Item(1).RowHeight = 15.65
$ws.Rows.Item(2).RowHeight = 44

# --- Sheet-wide default column width tweak ---
$ws.StandardWidth = 8.55859375

# --- Update the active selection to A2 ---
$ws.Range("A2").Select()
